# AutoCommit_12 декабря 2023 г. 16:55:14_SibNout2023
# Updates homework ("Дз") marks / variant ("Варианты") numbers for several
# students in the gradebook sheet, and moves the frozen-pane scroll position
# and active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits -------------------------------------------------

# Row 12 (student #8): variant filled in as 3
$ws.Range("K12").Value = 3

# Row 15 (student #11): variant filled in as 3
$ws.Range("K15").Value = 3

# Row 18 (student #14): variant corrected from 3 to 1
$ws.Range("K18").Value = 1

# Row 22 (student #18): extra column R filled in
$ws.Range("R22").Value = 5

# Row 25 (student #21): variant filled in as 2
$ws.Range("K25").Value = 2

# Row 26 (student #22): extra column R filled in
$ws.Range("R26").Value = 5

# Row 31 (student #27): homework marks filled in (Дз1, Дз2, Дз3, Дз4, Дз6)
$ws.Range("C31").Value = 5
$ws.Range("D31").Value = 5
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 5
$ws.Range("H31").Value = 5

# Row 32 (student #28): variant filled in as 3
$ws.Range("K32").Value = 3

# --- View state ---------------------------------------------------------
# Move the frozen-pane scroll position so row 5 / column C is the first
# visible (unfrozen) row/column, and select R9 in the bottom-right pane.

$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 3
$null = $ws.Range("R9").Select()
